$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows of market-price data (rows 2-6) get reordered: the two most recent
# records (previously rows 5-6, dated 45072) move to the top of the block
# (rows 2-3), and the three older records (previously rows 2-4) shift down to
# become rows 4-6. Row 7 is untouched. Only the columns that actually differ
# between records (D, L, M, N, O, P, R, S) show up in the diff; the rest of
# each row's columns are identical across all the affected rows, so setting
# just these columns reproduces the full reordering.

$ws.Range("D2").Value = 45072
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("R2").Value = "Provincia de Chacabuco"
$ws.Range("S2").Value = 889

$ws.Range("D3").Value = 45072
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 17000
$ws.Range("O3").Value = 17000
$ws.Range("P3").Value = 17000
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 944

$ws.Range("D4").Value = 44252
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 13500
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 750

$ws.Range("D5").Value = 44250
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 806

$ws.Range("D6").Value = 44253
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 160
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 806
